# Refresh the cryptos price/volume table (Price = column D, Volume(1h) = column E).
# Some "Price" strings look numeric (e.g. "210.66"), so Excel's COM auto-type-detection
# would silently turn them into Number cells. We force those specific cells to Text via
# NumberFormat "@" before assigning, then restore the "Normal" style so no stray
# formatting is left behind on cells whose value doesn't need the text coercion.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.089.37"
$ws.Range("D3").Value = "1.564.74"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.491"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.91"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.99%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0862"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("D12").Value = "1.787.46"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "1.565.01"
$ws.Range("E13").Value = "  +1.41%  "
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("D16").Value = "27.089.86"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.25%  "
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +1.45%  "
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("E30").Value = "  +4.49%  "
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("D34").Value = "1.442.31"
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.08"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("E36").Value = "  -0.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.92%  "
$ws.Range("E38").Value = "  +0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.530"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.59%  "
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("E43").Value = "  +1.09%  "
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "1.701.53"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.49%  "
$ws.Range("E49").Value = "  +3.93%  "
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0955"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.43%  "
